$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 195.2
$ws.Range("I42").Value = 60.11111
$ws.Range("K42").Value = 180.33333
$ws.Range("M42").Value = 49.66667000000001
$ws.Range("H74").Value = 5285.143
$ws.Range("I74").Value = 4996
$ws.Range("J74").Value = 5333.3335
$ws.Range("K74").Value = 4996
$ws.Range("L74").Value = 5333.3335
$ws.Range("M74").Value = -4060
$ws.Range("N74").Value = -7205.3335
$ws.Range("H76").Value = 3241.5
$ws.Range("I76").Value = 3241.5
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3241.5
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -2926.5
$ws.Range("N76").ClearContents()
$ws.Range("H77").Value = 5285.143
$ws.Range("I77").Value = 4996
$ws.Range("J77").Value = 5333.3335
$ws.Range("K77").Value = 24980
$ws.Range("L77").Value = 26666.6675
$ws.Range("M77").Value = -20300
$ws.Range("N77").Value = -36026.6675
$ws.Range("H79").Value = 3241.5
$ws.Range("I79").Value = 3241.5
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3241.5
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2149.5
$ws.Range("N79").ClearContents()
$ws.Range("H92").Value = 509.2
$ws.Range("I92").Value = 496.76923
$ws.Range("J92").Value = 590
$ws.Range("K92").Value = 496.76923
$ws.Range("L92").Value = 590
$ws.Range("M92").Value = 751.23077
$ws.Range("N92").Value = -3086
$ws.Range("H113").Value = 5758.8667
$ws.Range("I113").Value = 6730.6665
$ws.Range("J113").Value = 5111
$ws.Range("K113").Value = 6730.6665
$ws.Range("L113").Value = 5111
$ws.Range("M113").Value = -3476.6665
$ws.Range("N113").Value = -11619
$ws.Range("H132").Value = 5019.6924
$ws.Range("I132").Value = 3497
$ws.Range("J132").Value = 5971.375
$ws.Range("K132").Value = 10491
$ws.Range("L132").Value = 17914.125
$ws.Range("M132").Value = -7961
$ws.Range("N132").Value = -22974.125
$ws.Range("H133").Value = 88888
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H135").Value = 35291.9
$ws.Range("I135").Value = 771.7308
$ws.Range("K135").Value = 6945.577200000001
$ws.Range("M135").Value = -4410.577200000001
$ws.Range("H137").Value = 93154.45
$ws.Range("I137").Value = 2766.8333
$ws.Range("K137").Value = 8300.499899999999
$ws.Range("M137").Value = -5750.499899999999
$ws.Range("H138").Value = 1874
$ws.Range("I138").Value = 1352.0476
$ws.Range("J138").Value = 2870.4546
$ws.Range("K138").Value = 4056.142800000001
$ws.Range("L138").Value = 8611.363799999999
$ws.Range("M138").Value = 1083.857199999999
$ws.Range("N138").Value = -18891.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 6852.4443
$ws.Range("I31").Value = 3894.125
$ws.Range("K31").Value = 3894.125
$ws.Range("M31").Value = -3600.125
$ws.Range("H32").Value = 35790.34
$ws.Range("I32").Value = 19303.291
$ws.Range("J32").Value = 262487.25
$ws.Range("K32").Value = 19303.291
$ws.Range("L32").Value = 262487.25
$ws.Range("M32").Value = -19016.291
$ws.Range("N32").Value = -263061.25
$ws.Range("H61").Value = 2547.5
$ws.Range("I61").Value = 2261.0715
$ws.Range("J61").Value = 3550
$ws.Range("K61").Value = 2261.0715
$ws.Range("L61").Value = 3550
$ws.Range("M61").Value = -2049.0715
$ws.Range("N61").Value = -3974
$ws.Range("H74").Value = 2480.75
$ws.Range("I74").Value = 2480.75
$ws.Range("K74").Value = 2480.75
$ws.Range("M74").Value = -1606.75
$ws.Range("H77").Value = 2480.75
$ws.Range("I77").Value = 2480.75
$ws.Range("K77").Value = 12403.75
$ws.Range("M77").Value = -8035.75
$ws.Range("H122").Value = 10952.863
$ws.Range("I122").Value = 12720.277
$ws.Range("K122").Value = 38160.831
$ws.Range("M122").Value = -35710.831
$ws.Range("H132").Value = 11716.392
$ws.Range("J132").Value = 4830.1665
$ws.Range("L132").Value = 14490.4995
$ws.Range("N132").Value = -19550.4995
$ws.Range("H136").Value = 2547.5
$ws.Range("I136").Value = 2261.0715
$ws.Range("J136").Value = 3550
$ws.Range("K136").Value = 6783.2145
$ws.Range("L136").Value = 10650
$ws.Range("M136").Value = -4233.2145
$ws.Range("N136").Value = -15750

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H105").Value = 2781080.2
$ws.Range("I105").Value = 3033632.8
$ws.Range("K105").Value = 3033632.8
$ws.Range("M105").Value = -3031885.8
$ws.Range("H134").Value = 1043.6364
$ws.Range("I134").Value = 1071.7368
$ws.Range("K134").Value = 3215.2104
$ws.Range("M134").Value = -680.2103999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2104.6191
$ws.Range("H34").Value = 2104.6191
$ws.Range("H43").Value = 28438
$ws.Range("J43").Value = 28438
$ws.Range("L43").Value = 28438
$ws.Range("N43").Value = -28806
$ws.Range("H101").Value = 28438
$ws.Range("J101").Value = 28438
$ws.Range("L101").Value = 28438
$ws.Range("N101").Value = -34928
$ws.Range("H122").Value = 4330.154
$ws.Range("I122").Value = 4165.778
$ws.Range("J122").Value = 4700
$ws.Range("K122").Value = 12497.334
$ws.Range("L122").Value = 14100
$ws.Range("M122").Value = -10047.334
$ws.Range("N122").Value = -19000
$ws.Range("H132").Value = 2659.2
$ws.Range("I132").Value = 2495.1538
$ws.Range("J132").Value = 3725.5
$ws.Range("K132").Value = 7485.4614
$ws.Range("L132").Value = 11176.5
$ws.Range("M132").Value = -4955.4614
$ws.Range("N132").Value = -16236.5
$ws.Range("H133").Value = 63331.668
$ws.Range("J133").Value = 63331.668
$ws.Range("L133").Value = 63331.668
$ws.Range("N133").Value = -68391.66800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 5311.75
$ws.Range("J64").Value = 5311.75
$ws.Range("L64").Value = 15935.25
$ws.Range("N64").Value = -16475.25
$ws.Range("H67").Value = 5311.75
$ws.Range("J67").Value = 5311.75
$ws.Range("L67").Value = 15935.25
$ws.Range("N67").Value = -17807.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 26563.75
$ws.Range("J62").Value = 26563.75
$ws.Range("L62").Value = 26563.75
$ws.Range("N62").Value = -27935.75
$ws.Range("H65").Value = 26563.75
$ws.Range("J65").Value = 26563.75
$ws.Range("L65").Value = 79691.25
$ws.Range("N65").Value = -86555.25
$ws.Range("H102").Value = 2717.2307
$ws.Range("I102").Value = 2502.6
$ws.Range("K102").Value = 2502.6
$ws.Range("M102").Value = -880.5999999999999
$ws.Range("H132").Value = 3553.3572
$ws.Range("I132").Value = 3043.6
$ws.Range("J132").Value = 4827.75
$ws.Range("K132").Value = 9130.799999999999
$ws.Range("L132").Value = 14483.25
$ws.Range("M132").Value = -6600.799999999999
$ws.Range("N132").Value = -19543.25
$ws.Range("H141").Value = 46715
$ws.Range("J141").Value = 46715
$ws.Range("L141").Value = 46715
$ws.Range("N141").Value = -57075

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 27794.928
$ws.Range("I100").Value = 5151
$ws.Range("K100").Value = 5151
$ws.Range("M100").Value = -4610
$ws.Range("H122").Value = 11640.846
$ws.Range("I122").Value = 15916.375
$ws.Range("J122").Value = 4800
$ws.Range("K122").Value = 47749.125
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = -45299.125
$ws.Range("N122").Value = -19300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 21749.5
$ws.Range("J74").Value = 21749.5
$ws.Range("L74").Value = 21749.5
$ws.Range("N74").Value = -23621.5
$ws.Range("H77").Value = 21749.5
$ws.Range("J77").Value = 21749.5
$ws.Range("L77").Value = 65248.5
$ws.Range("N77").Value = -74608.5
$ws.Range("H100").Value = 2491.25
$ws.Range("I100").Value = 3102.9092
$ws.Range("K100").Value = 6205.8184
$ws.Range("M100").Value = -5664.8184
$ws.Range("H122").Value = 2737.111
$ws.Range("I122").Value = 1932.25
$ws.Range("K122").Value = 5796.75
$ws.Range("M122").Value = -3346.75
$ws.Range("H136").Value = 11166.5
$ws.Range("I136").Value = 11399.8
$ws.Range("K136").Value = 34199.39999999999
$ws.Range("M136").Value = -31649.39999999999
